# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 14 (pushing the existing rows
# 14-18 down to 15-19) in the Pomelo / Terminal Hortofrutícola Agro
# Chillán sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above row 14; everything currently at row 14 and
# below shifts down by one (old 14 -> 15, ..., old 18 -> 19).
$ws.Rows(14).Insert()

# Populate the newly inserted row 14 with the new weekly record.
$ws.Cells.Item(14, 1).Value  = 7
$ws.Cells.Item(14, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(14, 3).Value  = "Ñuble"
$ws.Cells.Item(14, 4).Value  = 45212
$ws.Cells.Item(14, 5).Value  = 16
$ws.Cells.Item(14, 6).Value  = "Fruta"
$ws.Cells.Item(14, 7).Value  = 100102
$ws.Cells.Item(14, 8).Value  = "Cítricos"
$ws.Cells.Item(14, 9).Value  = 100102006
$ws.Cells.Item(14, 10).Value = "Pomelo"
$ws.Cells.Item(14, 11).Value = "Start Ruby"
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 40
$ws.Cells.Item(14, 14).Value = 17000
$ws.Cells.Item(14, 15).Value = 17000
$ws.Cells.Item(14, 16).Value = 17000
$ws.Cells.Item(14, 17).Value = "`$/caja 14 kilos granel"
$ws.Cells.Item(14, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(14, 19).Value = 1214
$ws.Cells.Item(14, 20).Value = 14
